# Generate Report for Handoff
# Adds a new localization-status row (file 7a2f0b06-...) above the
# existing "acfbde12-..." row on all three worksheets (Overview, zh-cn, de-de).

$wb = $excel.ActiveWorkbook

$githubBase = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/d85aed4c46662ff704e50dfd180d621fc7737182/e2e/"

# ---------------------------------------------------------------------------
# Sheet 1: "Overview"
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item(1)

# Push existing data row (row 2) down to row 3, keeping its formatting.
$wsOverview.Rows.Item(2).Insert()

# Fix up the stale hyperlink left pointing at the now-empty row 2.
$wsOverview.Range("B2:B3").Hyperlinks.Delete()

# New row 2: 7a2f0b06 file
$wsOverview.Range("A2").Value = "7a2f0b06-cb68-4a3a-9c58-6f449971e259.md"
$wsOverview.Range("C2").Value = ".md"
$wsOverview.Range("D2").Value = "'"
$wsOverview.Range("E2").Value = "Ready for handoff"
$wsOverview.Range("F2").Value = "Ready for handoff"
$wsOverview.Range("G2").Value = "'2016-09-07 17:01:16"
$wsOverview.Hyperlinks.Add($wsOverview.Range("B2"), ($githubBase + "7a2f0b06-cb68-4a3a-9c58-6f449971e259.md"), [System.Type]::Missing, [System.Type]::Missing, "e2e\7a2f0b06-cb68-4a3a-9c58-6f449971e259.md")

# Row 3: acfbde12 file (restore the hyperlink that used to sit on row 2)
$wsOverview.Hyperlinks.Add($wsOverview.Range("B3"), ($githubBase + "acfbde12-9567-475a-903a-8d728b069320.md"), [System.Type]::Missing, [System.Type]::Missing, "e2e\acfbde12-9567-475a-903a-8d728b069320.md")

$loOverview = $wsOverview.ListObjects.Item(1)
$loOverview.Resize($wsOverview.Range("A1:G3"))

# ---------------------------------------------------------------------------
# Sheet 2: "zh-cn"
# ---------------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item(2)

$wsZhCn.Rows.Item(2).Insert()
$wsZhCn.Range("A2:A3").Hyperlinks.Delete()

$wsZhCn.Range("A2").Value = "7a2f0b06-cb68-4a3a-9c58-6f449971e259.md"
$wsZhCn.Range("B2").Value = ".md"
$wsZhCn.Range("C2").Value = "Ready for handoff"
$wsZhCn.Range("D2").Value = "e2e"
$wsZhCn.Range("E2").Value = "ht"
$wsZhCn.Range("F2").Value = "'False"
$wsZhCn.Range("G2").Value = "7a2f0b06-cb68-4a3a-9c58-6f449971e259.9d07e703626588052f69160e70a2e62890191b56.zh-cn.xlf"
$wsZhCn.Range("H2").Value = "'2016-09-07 17:00:57"
$wsZhCn.Range("I2").Value = "'"
$wsZhCn.Range("J2").Value = "'"
$wsZhCn.Range("K2").Value = "'0001-01-01 00:00:00"
$wsZhCn.Range("L2").Value = "'"
$wsZhCn.Range("M2").Value = "'True"
$wsZhCn.Range("N2").Value = "'"
$wsZhCn.Range("O2").Value = "'False"
$wsZhCn.Range("P2").Value = "'"
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A2"), ($githubBase + "7a2f0b06-cb68-4a3a-9c58-6f449971e259.md"), [System.Type]::Missing, [System.Type]::Missing, "7a2f0b06-cb68-4a3a-9c58-6f449971e259.md")

$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A3"), ($githubBase + "acfbde12-9567-475a-903a-8d728b069320.md"), [System.Type]::Missing, [System.Type]::Missing, "acfbde12-9567-475a-903a-8d728b069320.md")

$loZhCn = $wsZhCn.ListObjects.Item(1)
$loZhCn.Resize($wsZhCn.Range("A1:P3"))

# ---------------------------------------------------------------------------
# Sheet 3: "de-de"
# ---------------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item(3)

$wsDeDe.Rows.Item(2).Insert()
$wsDeDe.Range("A2:A3").Hyperlinks.Delete()

$wsDeDe.Range("A2").Value = "7a2f0b06-cb68-4a3a-9c58-6f449971e259.md"
$wsDeDe.Range("B2").Value = ".md"
$wsDeDe.Range("C2").Value = "Ready for handoff"
$wsDeDe.Range("D2").Value = "e2e"
$wsDeDe.Range("E2").Value = "ht"
$wsDeDe.Range("F2").Value = "'False"
$wsDeDe.Range("G2").Value = "7a2f0b06-cb68-4a3a-9c58-6f449971e259.9d07e703626588052f69160e70a2e62890191b56.de-de.xlf"
$wsDeDe.Range("H2").Value = "'2016-09-07 17:01:16"
$wsDeDe.Range("I2").Value = "'"
$wsDeDe.Range("J2").Value = "'"
$wsDeDe.Range("K2").Value = "'0001-01-01 00:00:00"
$wsDeDe.Range("L2").Value = "'"
$wsDeDe.Range("M2").Value = "'True"
$wsDeDe.Range("N2").Value = "'"
$wsDeDe.Range("O2").Value = "'False"
$wsDeDe.Range("P2").Value = "'"
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A2"), ($githubBase + "7a2f0b06-cb68-4a3a-9c58-6f449971e259.md"), [System.Type]::Missing, [System.Type]::Missing, "7a2f0b06-cb68-4a3a-9c58-6f449971e259.md")

$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A3"), ($githubBase + "acfbde12-9567-475a-903a-8d728b069320.md"), [System.Type]::Missing, [System.Type]::Missing, "acfbde12-9567-475a-903a-8d728b069320.md")

$loDeDe = $wsDeDe.ListObjects.Item(1)
$loDeDe.Resize($wsDeDe.Range("A1:P3"))
